$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
Write-Host ($win.SheetViews | Get-Member | Out-String)
